# Add a new worksheet named "outputs" after the existing sheet(s) and
# populate it with the scoring output data (bug fix for rounding).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "outputs"

# Header row
$ws.Range("B1").Value = "پاسخ سیستم"
$ws.Range("C1").Value = "پاسخ دستی کارشناس"

# Raw / t_score rows
$ws.Range("A2").Value = "raw"
$ws.Range("B2").Value = 137

$ws.Range("A3").Value = "t_score"
$ws.Range("B3").Value = 35

$ws.Range("A4").Value = "t_score_summary"
$ws.Range("B4").Value = "عدم رضایت از روابط زناشویی همسران"

$ws.Range("A5").Value = "marital_communication"
$ws.Range("B5").Value = 16

$ws.Range("A6").Value = "personality_issues"
$ws.Range("B6").Value = 24

$ws.Range("A7").Value = "religious_orientation"
$ws.Range("B7").Value = 16

$ws.Range("A8").Value = "Conflict_resolution"
$ws.Range("B8").Value = 14

$ws.Range("A9").Value = "financial_management"
$ws.Range("B9").Value = 15

$ws.Range("A10").Value = "leisure_activities"
$ws.Range("B10").Value = 13

$ws.Range("A11").Value = "sexual_relationship"
$ws.Range("B11").Value = 14

$ws.Range("A12").Value = "children_&_marriage"
$ws.Range("B12").Value = 13

$ws.Range("A13").Value = "family_&_friends"
$ws.Range("B13").Value = 12

$ws.Range("A14").Value = "personality_issues_interpretation"
$ws.Range("B14").Value = "رضایت زیاد"

$ws.Range("A15").Value = "marital_communication_interpretation"
$ws.Range("B15").Value = "رضایت متوسط"

$ws.Range("A16").Value = "Conflict_resolution_interpretation"
$ws.Range("B16").Value = "رضایت متوسط"

$ws.Range("A17").Value = "financial_management_interpretation"
$ws.Range("B17").Value = "رضایت متوسط"

$ws.Range("A18").Value = "leisure_activities_interpretation"
$ws.Range("B18").Value = "رضایت متوسط"

$ws.Range("A19").Value = "sexual_relationship_interpretation"
$ws.Range("B19").Value = "رضایت متوسط"

$ws.Range("A20").Value = "children_&_marriage_interpretation"
$ws.Range("B20").Value = "رضایت متوسط"

$ws.Range("A21").Value = "family_&_friends_interpretation"
$ws.Range("B21").Value = "عدم رضایت"

$ws.Range("A22").Value = "religious_orientation_interpretation"
$ws.Range("B22").Value = "رضایت متوسط"

$ws.Range("A1").Select()
